# Updates cryptos list prices/volumes (and reorders a few coin rows) to
# match the latest scrape, mirroring the upstream GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.147.05'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.853.61'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'237.74"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').Value = "'0.6860"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = "'0.07764"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.55%  '
$ws.Range('D9').Value = "'0.3040"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').Value = "'23.16"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').Value = "'0.08086"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = "'0.7209"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.193"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.807.63'
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('D15').Value = "'89.24"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '29.159.04'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = "'5.729"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.92%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = "'13.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.000007790"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = "'233.71"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.63%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '2.108.52'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = "'7.465"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.10%  '
$ws.Range('D25').Value = "'161.84"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = "'8.974"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').Value = "'0.1428"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').Value = "'1.953"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'1.409"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').Value = "'4.503"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('D32').Value = "'1.481"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.23%  '
$ws.Range('D33').Value = "'4.008"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = "'0.05205"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').Value = "'1.176"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('D36').Value = "'0.7037"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').Value = "'1.001"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = "'2.663"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = "'0.01850"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('D40').Value = "'2.690"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').Value = "'0.9333"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('D42').Value = '1.103.36'
$ws.Range('E42').Value = '  +5.42%  '
$ws.Range('D43').Value = "'0.4275"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').Value = "'5.899"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = "'70.16"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').Value = "'0.9999"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = "'102.41"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').Value = "'1.794"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = '2.003.87'
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').Value = "'9.141"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.44%  '
$ws.Range('D51').Value = "'6.988"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.90%  '
